# Auto-generated Excel COM-interop script to apply Odin_Profits scheduled-runner update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 3312.8096  # H5
$ws.Cells.Item(5, 9).Value = 991.53845  # I5
$ws.Cells.Item(5, 10).Value = 7084.875  # J5
$ws.Cells.Item(5, 11).Value = 991.53845  # K5
$ws.Cells.Item(5, 12).Value = 7084.875  # L5
$ws.Cells.Item(5, 13).Value = -876.53845  # M5
$ws.Cells.Item(5, 14).Value = -7314.875  # N5
$ws.Cells.Item(17, 8).Value = 450.61905  # H17
$ws.Cells.Item(17, 10).Value = 450.61905  # J17
$ws.Cells.Item(17, 12).Value = 1351.85715  # L17
$ws.Cells.Item(17, 14).Value = -1687.85715  # N17
$ws.Cells.Item(33, 8).Value = 426.3158  # H33
$ws.Cells.Item(33, 9).Value = 388.8889  # I33
$ws.Cells.Item(33, 11).Value = 388.8889  # K33
$ws.Cells.Item(33, 13).Value = -159.8889  # M33
$ws.Cells.Item(98, 8).Value = 6999.75  # H98
$ws.Cells.Item(98, 9).Value = 4333  # I98
$ws.Cells.Item(98, 11).Value = 4333  # K98
$ws.Cells.Item(98, 13).Value = -2835  # M98
$ws.Cells.Item(105, 8).Value = 27359.834  # H105
$ws.Cells.Item(105, 10).Value = 27359.834  # J105
$ws.Cells.Item(105, 12).Value = 27359.834  # L105
$ws.Cells.Item(105, 14).Value = -34347.834  # N105
$ws.Cells.Item(118, 8).Value = 1463.5714  # H118
$ws.Cells.Item(118, 9).Value = 1375.6666  # I118
$ws.Cells.Item(118, 11).Value = 4126.9998  # K118
$ws.Cells.Item(118, 13).Value = -2469.9998  # M118
$ws.Cells.Item(122, 8).Value = 6999.75  # H122
$ws.Cells.Item(122, 9).Value = 4333  # I122
$ws.Cells.Item(122, 11).Value = 12999  # K122
$ws.Cells.Item(122, 13).Value = -10549  # M122
$ws.Cells.Item(132, 8).Value = 10260.625  # H132
$ws.Cells.Item(132, 9).Value = 7103.875  # I132
$ws.Cells.Item(132, 11).Value = 21311.625  # K132
$ws.Cells.Item(132, 13).Value = -18781.625  # M132
$ws.Cells.Item(137, 8).Value = 9215.137000000001  # H137
$ws.Cells.Item(137, 10).Value = 16249.875  # J137
$ws.Cells.Item(137, 12).Value = 48749.625  # L137
$ws.Cells.Item(137, 14).Value = -53849.625  # N137
$ws.Cells.Item(138, 8).Value = 4464.273  # H138
$ws.Cells.Item(138, 10).Value = 4966.6787  # J138
$ws.Cells.Item(138, 12).Value = 14900.0361  # L138
$ws.Cells.Item(138, 14).Value = -25180.0361  # N138
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 600.5  # H4
$ws.Cells.Item(4, 9).Value = 600.5  # I4
$ws.Cells.Item(4, 11).Value = 600.5  # K4
$ws.Cells.Item(4, 13).Value = -484.5  # M4
$ws.Cells.Item(5, 8).Value = 708.4  # H5
$ws.Cells.Item(5, 9).Value = 260.5  # I5
$ws.Cells.Item(5, 11).Value = 260.5  # K5
$ws.Cells.Item(5, 13).Value = -148.5  # M5
$ws.Cells.Item(132, 8).Value = 26636  # H132
$ws.Cells.Item(132, 9).Value = 2052.0293  # I132
$ws.Cells.Item(132, 11).Value = 6156.0879  # K132
$ws.Cells.Item(132, 13).Value = -3626.0879  # M132
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 708.4  # H4
$ws.Cells.Item(4, 9).Value = 260.5  # I4
$ws.Cells.Item(4, 11).Value = 260.5  # K4
$ws.Cells.Item(4, 13).Value = -145.5  # M4
$ws.Cells.Item(22, 8).Value = 1667539.5  # H22
$ws.Cells.Item(22, 9).Value = 1047.4  # I22
$ws.Cells.Item(22, 11).Value = 1047.4  # K22
$ws.Cells.Item(22, 13).Value = -874.4000000000001  # M22
$ws.Cells.Item(64, 8).Value = 5701.75  # H64
$ws.Cells.Item(64, 9).Value = 1291  # I64
$ws.Cells.Item(64, 10).Value = 7172  # J64
$ws.Cells.Item(64, 11).Value = 1291  # K64
$ws.Cells.Item(64, 12).Value = 7172  # L64
$ws.Cells.Item(64, 13).Value = -1066  # M64
$ws.Cells.Item(64, 14).Value = -7622  # N64
$ws.Cells.Item(67, 8).Value = 5701.75  # H67
$ws.Cells.Item(67, 9).Value = 1291  # I67
$ws.Cells.Item(67, 10).Value = 7172  # J67
$ws.Cells.Item(67, 11).Value = 1291  # K67
$ws.Cells.Item(67, 12).Value = 7172  # L67
$ws.Cells.Item(67, 13).Value = -511  # M67
$ws.Cells.Item(67, 14).Value = -8732  # N67
$ws.Cells.Item(80, 8).Value = 15353.412  # H80
$ws.Cells.Item(80, 10).Value = 19225.23  # J80
$ws.Cells.Item(80, 12).Value = 19225.23  # L80
$ws.Cells.Item(80, 14).Value = -21221.23  # N80
$ws.Cells.Item(83, 8).Value = 15353.412  # H83
$ws.Cells.Item(83, 10).Value = 19225.23  # J83
$ws.Cells.Item(83, 12).Value = 96126.14999999999  # L83
$ws.Cells.Item(83, 14).Value = -106110.15  # N83
$ws.Cells.Item(86, 8).Value = 2460.7693  # H86
$ws.Cells.Item(86, 10).Value = 6206.3335  # J86
$ws.Cells.Item(86, 12).Value = 6206.3335  # L86
$ws.Cells.Item(86, 14).Value = -8452.333500000001  # N86
$ws.Cells.Item(89, 8).Value = 2460.7693  # H89
$ws.Cells.Item(89, 10).Value = 6206.3335  # J89
$ws.Cells.Item(89, 12).Value = 31031.6675  # L89
$ws.Cells.Item(89, 14).Value = -42263.6675  # N89
$ws.Cells.Item(107, 8).Value = 2957.353  # H107
$ws.Cells.Item(107, 9).Value = 3151.1538  # I107
$ws.Cells.Item(107, 10).Value = 2327.5  # J107
$ws.Cells.Item(107, 11).Value = 3151.1538  # K107
$ws.Cells.Item(107, 12).Value = 2327.5  # L107
$ws.Cells.Item(107, 13).Value = -1231.1538  # M107
$ws.Cells.Item(107, 14).Value = -6167.5  # N107
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5739.1665  # H31
$ws.Cells.Item(31, 9).Value = 991.1111  # I31
$ws.Cells.Item(31, 11).Value = 991.1111  # K31
$ws.Cells.Item(31, 13).Value = -696.1111  # M31
$ws.Cells.Item(34, 8).Value = 5739.1665  # H34
$ws.Cells.Item(34, 9).Value = 991.1111  # I34
$ws.Cells.Item(34, 11).Value = 991.1111  # K34
$ws.Cells.Item(34, 13).Value = -789.1111  # M34
$ws.Cells.Item(134, 8).Value = 565318.0600000001  # H134
$ws.Cells.Item(134, 9).Value = 1493332.5  # I134
$ws.Cells.Item(134, 11).Value = 4479997.5  # K134
$ws.Cells.Item(134, 13).Value = -4477462.5  # M134
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 15026.611  # H3
$ws.Cells.Item(3, 9).Value = 3407.9  # I3
$ws.Cells.Item(3, 10).Value = 29550  # J3
$ws.Cells.Item(3, 11).Value = 10223.7  # K3
$ws.Cells.Item(3, 12).Value = 88650  # L3
$ws.Cells.Item(3, 13).Value = -10111.7  # M3
$ws.Cells.Item(3, 14).Value = -88874  # N3
$ws.Cells.Item(108, 8).Value = 9193.714  # H108
$ws.Cells.Item(108, 9).Value = 871.2  # I108
$ws.Cells.Item(108, 10).Value = 30000  # J108
$ws.Cells.Item(108, 11).Value = 2613.6  # K108
$ws.Cells.Item(108, 12).Value = 90000  # L108
$ws.Cells.Item(108, 13).Value = 266.3999999999996  # M108
$ws.Cells.Item(108, 14).Value = -95760  # N108
$ws.Cells.Item(114, 8).Value = 665  # H114
$ws.Cells.Item(114, 9).Value = 699  # I114
$ws.Cells.Item(114, 10).Value = 495  # J114
$ws.Cells.Item(114, 11).Value = 2097  # K114
$ws.Cells.Item(114, 12).Value = 1485  # L114
$ws.Cells.Item(114, 13).Value = 1157  # M114
$ws.Cells.Item(114, 14).Value = -7993  # N114
$ws.Cells.Item(131, 8).Value = 45616604  # H131
$ws.Cells.Item(131, 10).Value = 22225340  # J131
$ws.Cells.Item(131, 12).Value = 66676020  # L131
$ws.Cells.Item(131, 14).Value = -66686100  # N131
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(31, 8).Value = 1438.375  # H31
$ws.Cells.Item(31, 9).Value = 1438.375  # I31
$ws.Cells.Item(31, 11).Value = 1438.375  # K31
$ws.Cells.Item(31, 13).Value = -1146.375  # M31
$ws.Cells.Item(37, 8).Value = 1438.375  # H37
$ws.Cells.Item(37, 9).Value = 1438.375  # I37
$ws.Cells.Item(37, 11).Value = 1438.375  # K37
$ws.Cells.Item(37, 13).Value = -1161.375  # M37
$ws.Cells.Item(70, 8).Value = 7874.72  # H70
$ws.Cells.Item(70, 9).Value = 7162.25  # I70
$ws.Cells.Item(70, 10).Value = 8532.385  # J70
$ws.Cells.Item(70, 11).Value = 7162.25  # K70
$ws.Cells.Item(70, 12).Value = 8532.385  # L70
$ws.Cells.Item(70, 13).Value = -6892.25  # M70
$ws.Cells.Item(70, 14).Value = -9072.385  # N70
$ws.Cells.Item(73, 8).Value = 7874.72  # H73
$ws.Cells.Item(73, 9).Value = 7162.25  # I73
$ws.Cells.Item(73, 10).Value = 8532.385  # J73
$ws.Cells.Item(73, 11).Value = 7162.25  # K73
$ws.Cells.Item(73, 12).Value = 8532.385  # L73
$ws.Cells.Item(73, 13).Value = -6226.25  # M73
$ws.Cells.Item(73, 14).Value = -10404.385  # N73
$ws.Cells.Item(80, 8).Value = 4746.9355  # H80
$ws.Cells.Item(80, 9).Value = 3155.7144  # I80
$ws.Cells.Item(80, 10).Value = 8088.5  # J80
$ws.Cells.Item(80, 11).Value = 3155.7144  # K80
$ws.Cells.Item(80, 12).Value = 8088.5  # L80
$ws.Cells.Item(80, 13).Value = -2157.7144  # M80
$ws.Cells.Item(80, 14).Value = -10084.5  # N80
$ws.Cells.Item(83, 8).Value = 4746.9355  # H83
$ws.Cells.Item(83, 9).Value = 3155.7144  # I83
$ws.Cells.Item(83, 10).Value = 8088.5  # J83
$ws.Cells.Item(83, 11).Value = 15778.572  # K83
$ws.Cells.Item(83, 12).Value = 40442.5  # L83
$ws.Cells.Item(83, 13).Value = -10786.572  # M83
$ws.Cells.Item(83, 14).Value = -50426.5  # N83
$ws.Cells.Item(101, 8).Value = 24903.25  # H101
$ws.Cells.Item(101, 10).Value = 24903.25  # J101
$ws.Cells.Item(101, 12).Value = 24903.25  # L101
$ws.Cells.Item(101, 14).Value = -31393.25  # N101
$ws.Cells.Item(102, 8).Value = 5846.6733  # H102
$ws.Cells.Item(102, 9).Value = 5019.8213  # I102
$ws.Cells.Item(102, 11).Value = 5019.8213  # K102
$ws.Cells.Item(102, 13).Value = -3397.8213  # M102
$ws.Cells.Item(122, 8).Value = 112348.3  # H122
$ws.Cells.Item(122, 9).Value = 180289.17  # I122
$ws.Cells.Item(122, 11).Value = 540867.51  # K122
$ws.Cells.Item(122, 13).Value = -538417.51  # M122
$ws.Cells.Item(126, 8).Value = 9964.117  # H126
$ws.Cells.Item(126, 9).Value = 4307  # I126
$ws.Cells.Item(126, 11).Value = 12921  # K126
$ws.Cells.Item(126, 13).Value = -10451  # M126
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 5165.6665  # H40
$ws.Cells.Item(40, 9).Value = 5497  # I40
$ws.Cells.Item(40, 11).Value = 5497  # K40
$ws.Cells.Item(40, 13).Value = -5361  # M40
$ws.Cells.Item(46, 8).Value = 1593.9166  # H46
$ws.Cells.Item(46, 9).Value = 978  # I46
$ws.Cells.Item(46, 11).Value = 978  # K46
$ws.Cells.Item(46, 13).Value = -790  # M46
$ws.Cells.Item(101, 8).Value = 46744.6  # H101
$ws.Cells.Item(101, 10).Value = 46744.6  # J101
$ws.Cells.Item(101, 12).Value = 46744.6  # L101
$ws.Cells.Item(101, 14).Value = -53234.6  # N101
$ws.Cells.Item(122, 8).Value = 4349.2856  # H122
$ws.Cells.Item(122, 9).Value = 3490  # I122
$ws.Cells.Item(122, 10).Value = 7500  # J122
$ws.Cells.Item(122, 11).Value = 10470  # K122
$ws.Cells.Item(122, 12).Value = 22500  # L122
$ws.Cells.Item(122, 13).Value = -8020  # M122
$ws.Cells.Item(122, 14).Value = -27400  # N122
$ws.Cells.Item(132, 8).Value = 3287.1667  # H132
$ws.Cells.Item(132, 9).Value = 2208.3333  # I132
$ws.Cells.Item(132, 11).Value = 6624.999899999999  # K132
$ws.Cells.Item(132, 13).Value = -4094.999899999999  # M132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 44999  # H4
$ws.Cells.Item(4, 9).Value = 44999  # I4
$ws.Cells.Item(4, 11).Value = 44999  # K4
$ws.Cells.Item(4, 13).Value = -44886  # M4
$ws.Cells.Item(81, 8).Value = 981.1667  # H81
$ws.Cells.Item(81, 9).Value = 981.1667  # I81
$ws.Cells.Item(81, 10).Value = 0  # J81
$ws.Cells.Item(81, 11).Value = 1962.3334  # K81
$ws.Cells.Item(81, 12).Value = 0  # L81
$ws.Cells.Item(81, 13).ClearContents()  # M81
$ws.Cells.Item(81, 14).Value = -901.3334  # N81
$ws.Cells.Item(84, 8).Value = 981.1667  # H84
$ws.Cells.Item(84, 9).Value = 981.1667  # I84
$ws.Cells.Item(84, 10).Value = 0  # J84
$ws.Cells.Item(84, 11).Value = 9811.666999999999  # K84
$ws.Cells.Item(84, 12).Value = 0  # L84
$ws.Cells.Item(84, 13).ClearContents()  # M84
$ws.Cells.Item(84, 14).Value = -4507.666999999999  # N84
$ws.Cells.Item(104, 8).Value = 31894.166  # H104
$ws.Cells.Item(104, 10).Value = 31894.166  # J104
$ws.Cells.Item(104, 12).Value = 31894.166  # L104
$ws.Cells.Item(104, 14).Value = -38882.166  # N104
$ws.Cells.Item(113, 8).Value = 3500.9048  # H113
$ws.Cells.Item(113, 9).Value = 2800.95  # I113
$ws.Cells.Item(113, 11).Value = 8402.849999999999  # K113
$ws.Cells.Item(113, 13).Value = -6232.849999999999  # M113
$ws.Cells.Item(122, 8).Value = 2994.7666  # H122
$ws.Cells.Item(122, 9).Value = 3030.9656  # I122
$ws.Cells.Item(122, 10).Value = 1945  # J122
$ws.Cells.Item(122, 11).Value = 9092.8968  # K122
$ws.Cells.Item(122, 12).Value = 5835  # L122
$ws.Cells.Item(122, 13).Value = -6642.8968  # M122
$ws.Cells.Item(122, 14).Value = -10735  # N122
